# Version 1.2: Add build AlarmList file func
#
# Adds a new "OPCDataTemplate" worksheet (after L1AlarmDataTemple) that holds
# template/instance strings used by the AlarmList file builder, and restores
# the various sheet selections / active-tab bookkeeping that Excel records
# when the workbook was last saved.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Add the new "OPCDataTemplate" sheet as the last (rightmost) tab.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet  = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "OPCDataTemplate"

# Header row
$newSheet.Range("A1").Value = "Template"
$newSheet.Range("B1").Value = "Instances"

# Template / instance pairs (column A written first, then column B, so the
# shared-string table picks up the same ordering as the authored workbook).
$newSheet.Range("A2").Value = "{0}_{1}_{2}_{3}_SIGNAL_{4}"
$newSheet.Range("A3").Value = "{0}_{1}_{2}_{3}_CM"
$newSheet.Range("A4").Value = "{0}_{1}_{2}_{3}_Hour"

$newSheet.Range("B2").Value = "CRISBELT_PLC01_TL35_BC045_SIGNAL_1"
$newSheet.Range("B3").Value = "CRISBELT_PLC07_EB03_BC019_CM"
$newSheet.Range("B4").Value = "CRISBELT_PLC07_EB03_BC019_Hour"

# Column widths sized (bestFit) to the widest entry in each column. The
# host's ColumnWidth setter quantises to 1/6-character steps, so feed it
# the pre-image that lands on the closest achievable step to the authored
# bestFit widths (25.6640625 / 41).
$newSheet.Columns.Item(1).ColumnWidth = 24.833333333333332
$newSheet.Columns.Item(2).ColumnWidth = 40.166666666666664

# Selection left on the sheet when the author saved the file.
$newSheet.Range("B4").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) Restore the remembered selections on the other sheets.
# ---------------------------------------------------------------------
$conveyor = $wb.Worksheets.Item("Conveyor")
$conveyor.Range("A2:D20").Select() | Out-Null

$baseList = $wb.Worksheets.Item("BaseList")
$baseList.Range("B2:B16").Select() | Out-Null

$testData = $wb.Worksheets.Item("TestDataTemplate")
$testData.Range("A7").Select() | Out-Null

# ---------------------------------------------------------------------
# 3) Make the new sheet the active / tab-selected sheet, matching the
#    bookmark the author left (activeTab points at OPCDataTemplate).
# ---------------------------------------------------------------------
$newSheet.Activate() | Out-Null
$newSheet.Range("B4").Select() | Out-Null
